$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the blank cell A2 by shifting the three keyword cells below it
# (A3:A5) up by one row (value + formatting together). This only affects
# rows 2-5; rows 6:9 (the trailing blank placeholder cells) are left as-is.
$ws.Cells.Item(3, 1).Copy($ws.Cells.Item(2, 1))
$ws.Cells.Item(4, 1).Copy($ws.Cells.Item(3, 1))
$ws.Cells.Item(5, 1).Copy($ws.Cells.Item(4, 1))

# The vacated row (old A5) becomes fully blank again.
$ws.Cells.Item(5, 1).Clear()

# Reflect the resulting selection on the sheet (A4 is selected after the edit).
$ws.Range("A4").Select()
